$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SpecimentsInfo")
$ws.Name = "SpecimenInfo"
